$d = $word.ActiveDocument

# 1) First paragraph: apply the "Body Text" style and swap its text.
$p1 = $d.Paragraphs(1)
$p1.Style = "Body Text"
$d.Content.Find.Execute("Что красное и пахнет как синяя краска?", $true, $false, $false, $false, $false, $true, 1, $false, "В каком году Эстония была провозглашена независимsм государством.", 2)

# 2) Second paragraph: swap its text (bookmark stays attached to this paragraph for now).
$d.Content.Find.Execute("Что зеленое и имеет колеса?", $true, $false, $false, $false, $false, $true, 1, $false, "Сколько Чудес Света существует в мире?", 2)

# 3) Insert a brand-new third paragraph after paragraph 2, moving the _GoBack
#    bookmark's trailing position to the end of the document, and set its text.
$p2 = $d.Paragraphs(2)
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs(3)
$p3.Range.Text = "Сколько Континентов на земле?"
